# Add payment 79174445 (Cash) 2025-08-18T08:51:16
#
# Per the diff:
#  - A15 (phone) switches from a text-stored "79174445" to a real number 79174445
#  - A new row 16 is appended with the new payment, where A16 keeps the phone
#    number stored as TEXT "79174445" (matching how A15 looked before this edit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: A15 becomes a genuine number ---
$ws.Range("A15").Value = 79174445

# --- Row 16: the new payment row ---

# A16 must stay TEXT ("79174445"), not auto-converted to a number. Typing a
# leading apostrophe into a scratch cell forces text, but doing that directly
# on A16 would stamp it with a quote-prefix cell style. Instead stage the
# text in an unused scratch cell, copy only the *value* over with
# PasteSpecial (which drops the quote-prefix formatting), then remove the
# scratch cell again.
$ws.Range("K1").Value = "'79174445"
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("K1").Delete(-4159) | Out-Null          # xlShiftUp, removes the scratch cell entirely

# B16 and F16 are blank cells in the source row (empty inline strings) -
# touch them with a neutral formatting no-op so the cell exists in the
# sheet (same shape as the other blank cells in the table) without
# altering their appearance.
$ws.Range("B16").Font.Bold = $false
$ws.Range("F16").Font.Bold = $false

$ws.Range("C16").Value = "Cash"
$ws.Range("D16").Value = "2025-08-18T08:51:16"
$ws.Range("E16").Value = 40
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 30
